# Generate Report for Archive
#
# 1) Shared-string text change: "Ready for handoff" -> "In Translation".
#    That string is referenced from Overview!E2, Overview!F2, zh-cn!C2 and
#    de-de!C2, so all four cells are rewritten to the new text.
# 2) Column-width shrink from 17.2159881591797 -> 13.4101848602295 chars on:
#       Overview  columns E (5) and F (6)
#       zh-cn     column  C (3)
#       de-de     column  C (3)
#    The host's ColumnWidth setter quantizes to an integer "MDW-6" pixel
#    grid (raw = (round(width*6)+5)/6), so the nearest representable value
#    to the target is 13.333333333333334 (round(width*6)=80) -- that's what
#    a ColumnWidth of 12.5 produces, comfortably inside the rounding bucket.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- text: "Ready for handoff" -> "In Translation" ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- column widths ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
